$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "Record"
$ws.Range("B56").Value = "RJ Record"
$ws.Range("C56").Value = "Defesa Civil"
$ws.Range("D56").Value = "2025-04-04T18:11"
$ws.Range("E56").Value = "Positivo"
$ws.Range("F56").Value = "Após previsão de fortes chuvas, Defesa Civil intensifica monitoramento. Entrevista com o coord. do Centro de Monitoramento de Desastres, Leandro Freitas. Previsão de temporais no final de semana. Últimos dias marcados por sol e calor. Nova frente fria chegando no Estado. Enquete com população. Algumas pessoas citaram pontos de alagamento. Estão sendo esperados 120 mm de chuva. Defesa Civil acompanha deslocamento da frente fria. Volume pode superar o esperado para o mês de abril. Defesa Civil incentiva a aderir ao sistema de alerta. Para se cadastrar, basta enviar o CEP para 40199. *matéria*"

$ws.Range("A57").Value = "Record"
$ws.Range("B57").Value = "RJ Record"
$ws.Range("C57").Value = "Saúde"
$ws.Range("D57").Value = "2025-04-04T18:39"
$ws.Range("E57").Value = "Positivo"
$ws.Range("F57").Value = "Vacinação contra o vírus influenza começa na próxima segunda-feira. Repórter *ao vivo*. Campanha será realizada em uma única etapa. De acordo com Secretaria de Saúde, vão ser divulgados os locais até o fim desta semana. Crianças de 5 anos, gestantes, idosos e, este ano, funcionários dos Correios e da área administrativa fazem parte do público alvo. Tema também foi abordado no Balanço Geral. "

$ws.Range("A58").Value = "Record"
$ws.Range("B58").Value = "RJ Record"
$ws.Range("C58").Value = "CCZ"
$ws.Range("D58").Value = "2025-04-04T18:56"
$ws.Range("E58").Value = "Positivo"
$ws.Range("F58").Value = "Animais são resgatados após denúncia de maus tratos. Caso foi parar na polícia civil, após denúncia anônima. Ação conjunta entre a Polícia Civil e CCZ foi realizada. Caso aconteceu na Rua Augusto Bessa, no Turfe. Dona da casa teria ido para Farol, sem data para voltar. Animais estavam acorrentados há cerca de 15 dias. Na residência, foram encontrados um gato e uma cachorra. Água estava com larvas de mosquito. Animais foram levados para o canil do CCZ. Entrevista com veterinário do CCZ, José Leonardo. Exibido vídeo dos animais acorrentados. Esses não foram encontrados. Eles estavam em ambiente insalubre. Ele também convidou as pessoas a aderir à adoção responsável.  Maus tratos é crime. Abril Laranja é o mês de conscientização e combate aos maus-tratos a animais. Tema também foi abordado no Balanço Geral. *matéria*"
